{"js": "// Helper: wrap a run-level OOXML fragment into a single <w:p> so it can be\n// inserted via insertOoxml(..., Word.InsertLocation.replace) \u2014 this keeps\n// all the <w:t>/<w:br/> children inside ONE <w:r>, matching the canonical\n// OOXML produced by the original authoring tool (instead of Office.js's\n// default of one <w:r> per insertText/insertBreak call).\nfunction wrapParagraph(innerRunXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' + innerRunXml + '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nfunction escXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst paras = body.paragraphs.items;\n\n// Paragraph 1 (0-based index 1): was 3 lines about Organism K/F/J/G feeding\n// relationships, now keeps only the first line text (also fixing\n// \"feeds-on\" -> \"feeds on\").\nparas[1].insertOoxml(\n  wrapParagraph(\n    \"<w:r><w:t>\" + escXml(\". Organism K feeds on Organism. J.\") + \"</w:t></w:r>\"\n  ),\n  Word.InsertLocation.replace\n);\n\n// Paragraph 2 (0-based index 2): was the first inline picture (7.jpg);\n// it is replaced by the two remaining \"Organism\" lines that used to live\n// in paragraph 1, now as a two-line text paragraph (with <w:br/> between).\nparas[2].insertOoxml(\n  wrapParagraph(\n    \"<w:r><w:t>\" + escXml(\"\\u2018 Organism F and J feed on Organism H.\") + \"</w:t>\" +\n    \"<w:br/>\" +\n    \"<w:t>\" + escXml(\"e| Organism G.feeds on Organisms F and H:\") + \"</w:t></w:r>\"\n  ),\n  Word.InsertLocation.replace\n);\n\n// Paragraph 3 (0-based index 3): plain typo/OCR-cleanup fixes.\nparas[3].insertOoxml(\n  wrapParagraph(\n    \"<w:r><w:t>\" +\n      escXml(\"6. -. \\u2018Fhe picture below shows the structural adaptations of two piants.\") +\n      \"</w:t></w:r>\"\n  ),\n  Word.InsertLocation.replace\n);\n\n// Paragraph 4 (0-based index 4): was the second inline picture (6.jpg);\n// the whole paragraph is removed.\nparas[4].delete();\n\n// Paragraph 5 (0-based index 5): only a single word typo fix\n// (\"folowing\" -> \"following\"); everything else (breaks, other runs of\n// text) stays exactly as-is, so do a targeted search & replace instead of\n// rebuilding the whole paragraph.\nconst search = body.search(\"folowing\", { matchCase: true, matchWholeWord: false });\nsearch.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < search.items.length; i++) {\n  search.items[i].insertText(\"following\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 2: was 3 lines (Organism K / F+J / G feeding relationships).\n# Trim it down to just the first line, fixing \"feeds-on\" -> \"feeds on\".\n$p2 = $d.Paragraphs.Item(2)\n$p2.Range.Text = \". Organism K feeds on Organism. J.\"\n\n# Paragraph 3: was the first inline picture (7.jpg). Replace it with the\n# two lines that used to be part of paragraph 2 (now as plain text,\n# separated by a manual line break, char 11 == <w:br/>).\n$p3 = $d.Paragraphs.Item(3)\n$p3.Range.Text = [char]0x2018 + \" Organism F and J feed on Organism H.\" + [char]11 + \"e| Organism G.feeds on Organisms F and H:\"\n\n# Paragraph 4: plain typo / OCR-cleanup fixes.\n$p4 = $d.Paragraphs.Item(4)\n$p4.Range.Text = \"6. -. \" + [char]0x2018 + \"Fhe picture below shows the structural adaptations of two piants.\"\n\n# Paragraph 5: was the second inline picture (6.jpg). Delete the whole\n# paragraph (text + picture + its paragraph mark).\n$p5 = $d.Paragraphs.Item(5)\n$rng = $d.Range($p5.Range.Start, $p5.Range.End)\n$rng.Delete()\n\n# Paragraph 6 (now last paragraph): single word typo fix only\n# (\"folowing\" -> \"following\"); leave every other run / break untouched.\n$find = $d.Content\n$find.Find.Execute(\"folowing\", $false, $false, $false, $false, $false, $true, 1, $false, \"following\", 2)\n"}
